$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the recruitment band elevation criteria values
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 150

# Move active selection to C25 (matches the saved worksheet view state)
$ws.Range("C25").Select()
